# Ersatt "hårdkodade" värden under avsnitt 4
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Månad" (month number) column B from 10 -> 11 for data rows 2-13
# (stored as text in the sheet, not a number, so force text format to
# avoid Excel auto-converting the numeric-looking string into a number)
$rngManad = $ws.Range("B2:B13")
$rngManad.NumberFormat = "@"
$rngManad.Value = "11"
$rngManad.Style = "Normal"

# Update "månad_namn" (month name) column J from oktober -> november for data rows 2-13
$ws.Range("J2:J13").Value = "november"

# Update "Ohälsotalet" values in column G
$ws.Cells.Item(2, 7).Value = 22.5
$ws.Cells.Item(3, 7).Value = 28.3
$ws.Cells.Item(4, 7).Value = 26.7
$ws.Cells.Item(5, 7).Value = 28.2
$ws.Cells.Item(6, 7).Value = 27.2
$ws.Cells.Item(7, 7).Value = 35.2
$ws.Cells.Item(8, 7).Value = 32.6
$ws.Cells.Item(9, 7).Value = 34.5
$ws.Cells.Item(10, 7).Value = 18.1
$ws.Cells.Item(11, 7).Value = 21.9
$ws.Cells.Item(12, 7).Value = 21.2
$ws.Cells.Item(13, 7).Value = 22.3
